$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.37"
$ws.Range("E2").Value = "'-2.84%"
$ws.Range("G2").Value = "'19"
$ws.Range("D3").Value = "'40.76"
$ws.Range("E3").Value = "'-2.55%"
$ws.Range("G3").Value = "'19"
$ws.Range("D4").Value = "'5.055"
$ws.Range("E4").Value = "'-2.55%"
$ws.Range("G4").Value = "'19"
$ws.Range("D5").Value = "'0.07619"
$ws.Range("E5").Value = "'-5.48%"
$ws.Range("G5").Value = "'19"
$ws.Range("D6").Value = "'4.245"
$ws.Range("E6").Value = "'-2.85%"
$ws.Range("G6").Value = "'19"
$ws.Range("D7").Value = "'1.592"
$ws.Range("E7").Value = "'-8.58%"
$ws.Range("G7").Value = "'19"
$ws.Range("D8").Value = "'0.9065"
$ws.Range("E8").Value = "'-2.34%"
$ws.Range("G8").Value = "'19"
$ws.Range("E9").Value = "'-12.66%"
$ws.Range("G9").Value = "'19"
$ws.Range("D10").Value = "'0.1766"
$ws.Range("E10").Value = "'-4.90%"
$ws.Range("G10").Value = "'19"
$ws.Range("D11").Value = "'0.09089"
$ws.Range("E11").Value = "'-2.59%"
$ws.Range("G11").Value = "'19"
$ws.Range("D12").Value = "'0.04383"
$ws.Range("E12").Value = "'-3.95%"
$ws.Range("G12").Value = "'19"
$ws.Range("E13").Value = "'-0.05%"
$ws.Range("G13").Value = "'19"
$ws.Range("D14").Value = "'0.001251"
$ws.Range("E14").Value = "'-2.11%"
$ws.Range("G14").Value = "'19"
$ws.Range("D15").Value = "'0.005799"
$ws.Range("E15").Value = "'-2.53%"
$ws.Range("G15").Value = "'19"
$ws.Range("D16").Value = "'3.367"
$ws.Range("E16").Value = "'0.36%"
$ws.Range("G16").Value = "'19"
$ws.Range("D17").Value = "'2.440"
$ws.Range("E17").Value = "'-5.06%"
$ws.Range("G17").Value = "'19"
$ws.Range("D18").Value = "'0.3299"
$ws.Range("E18").Value = "'-2.56%"
$ws.Range("G18").Value = "'19"
$ws.Range("D19").Value = "'6.839"
$ws.Range("E19").Value = "'-7.16%"
$ws.Range("G19").Value = "'19"
$ws.Range("D20").Value = "'0.1339"
$ws.Range("E20").Value = "'-3.15%"
$ws.Range("G20").Value = "'19"
$ws.Range("D21").Value = "'0.2843"
$ws.Range("E21").Value = "'9.15%"
$ws.Range("G21").Value = "'19"
$ws.Range("D22").Value = "'0.04155"
$ws.Range("E22").Value = "'-0.40%"
$ws.Range("G22").Value = "'19"
$ws.Range("D23").Value = "'0.001213"
$ws.Range("E23").Value = "'-2.86%"
$ws.Range("G23").Value = "'19"
$ws.Range("D24").Value = "'0.004008"
$ws.Range("E24").Value = "'-6.30%"
$ws.Range("G24").Value = "'19"
$ws.Range("E25").Value = "'6.22%"
$ws.Range("G25").Value = "'19"
$ws.Range("D26").Value = "'0.0003008"
$ws.Range("E26").Value = "'0.68%"
$ws.Range("G26").Value = "'19"
$ws.Range("G27").Value = "'19"
$ws.Range("G28").Value = "'19"
$ws.Range("G29").Value = "'19"
$ws.Range("G30").Value = "'19"
$ws.Range("G31").Value = "'19"
$ws.Range("G32").Value = "'19"
$ws.Range("G33").Value = "'19"
$ws.Range("G34").Value = "'19"
$ws.Range("G35").Value = "'19"
$ws.Range("G36").Value = "'19"
$ws.Range("G37").Value = "'19"
$ws.Range("D38").Value = "'0.02415"
$ws.Range("E38").Value = "'-5.61%"
$ws.Range("G38").Value = "'19"
$ws.Range("D39").Value = "'0.05144"
$ws.Range("E39").Value = "'-5.08%"
$ws.Range("G39").Value = "'19"
$ws.Range("D40").Value = "'0.007833"
$ws.Range("E40").Value = "'-2.83%"
$ws.Range("G40").Value = "'19"
$ws.Range("D41").Value = "'0.1308"
$ws.Range("E41").Value = "'-5.88%"
$ws.Range("G41").Value = "'19"
$ws.Range("D42").Value = "'0.007075"
$ws.Range("E42").Value = "'-6.67%"
$ws.Range("G42").Value = "'19"
$ws.Range("D43").Value = "'0.001949"
$ws.Range("E43").Value = "'-1.92%"
$ws.Range("G43").Value = "'19"
$ws.Range("D44").Value = "'0.008036"
$ws.Range("E44").Value = "'-3.45%"
$ws.Range("G44").Value = "'19"
$ws.Range("D45").Value = "'0.3052"
$ws.Range("E45").Value = "'-2.73%"
$ws.Range("G45").Value = "'19"
$ws.Range("D46").Value = "'0.00006366"
$ws.Range("E46").Value = "'-6.01%"
$ws.Range("G46").Value = "'19"
$ws.Range("E47").Value = "'-0.26%"
$ws.Range("G47").Value = "'19"
$ws.Range("E48").Value = "'-26.97%"
$ws.Range("G48").Value = "'19"
$ws.Range("D49").Value = "'0.005716"
$ws.Range("E49").Value = "'68.35%"
$ws.Range("G49").Value = "'19"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.26%"
$ws.Range("G50").Value = "'19"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.26%"
$ws.Range("G51").Value = "'19"
